$wb = $excel.ActiveWorkbook

# Move the "BPH" sheet so it is positioned before the "field" sheet.
$bph = $wb.Worksheets.Item("BPH")
$field = $wb.Worksheets.Item("field")
$bph.Move($field)

# Update the simulation parameters on the (now reordered) "BPH" sheet:
# stop spraying pesticide during weeks corresponding to rows 4-8 and 11.
$bphWs = $wb.Worksheets.Item("BPH")
$bphWs.Range("C4").Value = 0
$bphWs.Range("C5").Value = 0
$bphWs.Range("C6").Value = 0
$bphWs.Range("C7").Value = 0
$bphWs.Range("C8").Value = 0
$bphWs.Range("C11").Value = 0

# Update the calcul sheet's initial pest-count parameter from 300 to 100.
$calcWs = $wb.Worksheets.Item("calcul")
$calcWs.Range("N2").Formula = "=100"
